# Update bitcoin_buys.xlsx after running on 2025-05-07
# Appends a new row of data (row 13) to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 13

# Column A: date value stored as plain text (matches the style of the other
# "MM/DD/YYYY" rows already in the sheet, e.g. rows 10 and 12).
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "05/07/2025"
$dateCell.Style = "Normal"

# Columns B-D: numeric values.
$ws.Cells.Item($row, 2).Value = 0.0005152999999999998
$ws.Cells.Item($row, 3).Value = 97030.85581214831
$ws.Cells.Item($row, 4).Value = 50
